$d = $word.ActiveDocument

# Change 1: Overview paragraph 1 - rephrase experience sentence
$d.Paragraphs(4).Range.Text = "I am an engineering leader who is passionate about building great teams with a focus on autonomy and accountability. Over the last seven years, I have led and mentored engineers and managers. I bring an extensive background building cloud-based distributed systems with a strong emphasis on DevOps culture"

# Change 2: Overview paragraph 2 - rephrase Nike Tech Talks / conferences sentence
$d.Paragraphs(5).Range.Text = "Since 2015, I have hosted the Nike Tech Talks, a highly successful event series showcasing world class subject matter experts from across the industry. I am also the organizer of Nike’s open-source program, and a champion of building a strong and healthy engineering culture. I have organized two conferences (Pacific Northwest Scala) as well as founding and running multiple user group communities. And I have been an invited speaker at major industry conferences and special interest group events."

# Change 3: Accomplishments - $250m -> $250M
$d.Paragraphs(7).Range.Text = "Led the design and build out of a communications platform responsible for over `$250M in attributable revenue."

# Change 4: Accomplishments - fitness tracking platform sentence gains a clause
$d.Paragraphs(9).Range.Text = "Led design, development, and ongoing support of a fitness activity tracking platform that scaled to handle Nike’s global consumer ecosystem."

# Change 5: Experience / Nike Sr. Engineering Manager description
$d.Paragraphs(18).Range.Text = "Leading a multiple teams across backend and frontend to provide a unified, real-time communications platform. Established the early groundwork of the platform and helped define and drive the larger platform strategy. This is still a very active project with a business value of at least `$250M in attributable revenue."
